# Scheduled runner refresh: re-pulls Universalis market prices (currentAveragePrice*
# and the derived LevePrice*/LeveProfit* columns, H:N) for the affected leves across
# the Shiva_Profits crafting-job sheets. Item/leve identity columns (A:G) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: Enchanted Silver Ink
$ws.Range("H28").Value = 831.74286
$ws.Range("I28").Value = 765.5517
$ws.Range("J28").Value = 1151.6666
$ws.Range("K28").Value = 765.5517
$ws.Range("L28").Value = 1151.6666
$ws.Range("M28").Value = -280.5517
$ws.Range("N28").Value = -2121.6666
# Row 62: Enchanted Mythrite Ink
$ws.Range("H62").Value = 5898.125
$ws.Range("I62").Value = 5924.75
$ws.Range("J62").Value = 5871.5
$ws.Range("K62").Value = 5924.75
$ws.Range("L62").Value = 5871.5
$ws.Range("M62").Value = -5300.75
$ws.Range("N62").Value = -7119.5
# Row 65: Enchanted Mythrite Ink
$ws.Range("H65").Value = 5898.125
$ws.Range("I65").Value = 5924.75
$ws.Range("J65").Value = 5871.5
$ws.Range("K65").Value = 29623.75
$ws.Range("L65").Value = 29357.5
$ws.Range("M65").Value = -26503.75
$ws.Range("N65").Value = -35597.5
# Row 100: Beetle Glue
$ws.Range("H100").Value = 1998.238
$ws.Range("I100").Value = 1988.5555
$ws.Range("K100").Value = 1988.5555
$ws.Range("M100").Value = -1447.5555
# Row 101: Cunning Craftsman's Tea
$ws.Range("H101").Value = 403.27274
$ws.Range("I101").Value = 403.27274
$ws.Range("K101").Value = 1209.81822
$ws.Range("M101").Value = 412.1817799999999
# Row 106: Enchanted Palladium Ink
$ws.Range("H106").Value = 3714.2856
$ws.Range("I106").Value = 1999
$ws.Range("K106").Value = 1999
$ws.Range("M106").Value = -1368
# Row 107: Enchanted Truegold Ink
$ws.Range("H107").Value = 52633664
$ws.Range("I107").Value = 66667544
$ws.Range("K107").Value = 66667544
$ws.Range("M107").Value = -66665624
# Row 115: Competent Craftsman's Syrup
$ws.Range("H115").Value = 125006680
$ws.Range("I115").Value = 125006680
$ws.Range("K115").Value = 375020040
$ws.Range("M115").Value = -375018473
# Row 125: Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 2244
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2244
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 20196
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -25116
# Row 127: Competent Craftsman's Draught
$ws.Range("H127").Value = 5883.1177
$ws.Range("I127").Value = 1507.5
$ws.Range("K127").Value = 4522.5
$ws.Range("M127").Value = 437.5
# Row 138: Cunning Craftsman's Tisane
$ws.Range("H138").Value = 37044710
$ws.Range("J138").Value = 10479.066
$ws.Range("L138").Value = 31437.198
$ws.Range("N138").Value = -41717.198

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Bronze Ingot
$ws.Range("H2").Value = 1485.8572
$ws.Range("I2").Value = 1487.9333
$ws.Range("K2").Value = 1487.9333
$ws.Range("M2").Value = -1374.9333
# Row 32: Steel Ingot
$ws.Range("H32").Value = 1470.44
$ws.Range("I32").Value = 1191.6022
$ws.Range("J32").Value = 5175
$ws.Range("K32").Value = 1191.6022
$ws.Range("L32").Value = 5175
$ws.Range("M32").Value = -904.6022
$ws.Range("N32").Value = -5749
# Row 102: Tama-hagane Ingot
$ws.Range("H102").Value = 7146.0713
$ws.Range("I102").Value = 6707.375
$ws.Range("K102").Value = 6707.375
$ws.Range("M102").Value = -5085.375
# Row 110: Deepgold Ingot
$ws.Range("H110").Value = 1703.3549
$ws.Range("I110").Value = 1553.8334
$ws.Range("J110").Value = 2216
$ws.Range("K110").Value = 1553.8334
$ws.Range("L110").Value = 2216
$ws.Range("M110").Value = 491.1666
$ws.Range("N110").Value = -6306
# Row 116: Titanbronze Ingot
$ws.Range("H116").Value = 1485.8572
$ws.Range("I116").Value = 1487.9333
$ws.Range("K116").Value = 1487.9333
$ws.Range("M116").Value = 806.0667000000001
# Row 122: High Durium Nugget
$ws.Range("H122").Value = 7790.947
$ws.Range("I122").Value = 6498.9
$ws.Range("J122").Value = 8260.781999999999
$ws.Range("K122").Value = 19496.7
$ws.Range("L122").Value = 24782.346
$ws.Range("M122").Value = -17046.7
$ws.Range("N122").Value = -29682.346
# Row 132: Mountain Chromite Ingot
$ws.Range("H132").Value = 3315
$ws.Range("I132").Value = 3454.6155
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 10363.8465
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -7833.8465
$ws.Range("N132").Value = -9560

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Bronze Ingot
$ws.Range("H3").Value = 1485.8572
$ws.Range("I3").Value = 1487.9333
$ws.Range("K3").Value = 1487.9333
$ws.Range("M3").Value = -1373.9333
# Row 105: Molybdenum Ingot
$ws.Range("H105").Value = 1790.5151
$ws.Range("I105").Value = 1593.4615
$ws.Range("K105").Value = 1593.4615
$ws.Range("M105").Value = 153.5385000000001
# Row 107: Deepgold Nugget
$ws.Range("H107").Value = 3750.3442
$ws.Range("I107").Value = 3309.6316
$ws.Range("J107").Value = 4478.478
$ws.Range("K107").Value = 3309.6316
$ws.Range("L107").Value = 4478.478
$ws.Range("M107").Value = -1389.6316
$ws.Range("N107").Value = -8318.477999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 94: Beech Lumber
$ws.Range("H94").Value = 3668.8333
$ws.Range("J94").Value = 3606.875
$ws.Range("L94").Value = 3606.875
$ws.Range("N94").Value = -4508.875
# Row 107: White Oak Lumber
$ws.Range("H107").Value = 1612.1482
$ws.Range("I107").Value = 791.73334
$ws.Range("K107").Value = 791.73334
$ws.Range("M107").Value = 1128.26666

$ws = $wb.Worksheets.Item("CUL")
# Row 11: Orange Juice
$ws.Range("H11").Value = 1378.4736
$ws.Range("I11").Value = 837
$ws.Range("J11").Value = 4266.3335
$ws.Range("K11").Value = 2511
$ws.Range("L11").Value = 12799.0005
$ws.Range("M11").Value = -2371
$ws.Range("N11").Value = -13079.0005
# Row 38: Dark Vinegar
$ws.Range("H38").Value = 84.77273
$ws.Range("I38").Value = 72.64706
$ws.Range("J38").Value = 126
$ws.Range("K38").Value = 217.94118
$ws.Range("L38").Value = 378
$ws.Range("M38").Value = 129.05882
$ws.Range("N38").Value = -1072
# Row 46: Acorn Cookie
$ws.Range("H46").Value = 140240450
$ws.Range("I46").Value = 217.1579
$ws.Range("J46").Value = 288271800
$ws.Range("K46").Value = 651.4737
$ws.Range("L46").Value = 864815400
$ws.Range("M46").Value = -560.4737
$ws.Range("N46").Value = -864815582
# Row 68: Fermented Butter
$ws.Range("H68").Value = 3377.3333
$ws.Range("I68").Value = 1196
$ws.Range("K68").Value = 3588
$ws.Range("M68").Value = -2777
# Row 71: Fermented Butter
$ws.Range("H71").Value = 3377.3333
$ws.Range("I71").Value = 1196
$ws.Range("K71").Value = 10764
$ws.Range("M71").Value = -6708

$ws = $wb.Worksheets.Item("GSM")
# Row 46: Fire Brand
$ws.Range("H46").Value = 57588.4
$ws.Range("I46").Value = 33314
$ws.Range("K46").Value = 33314
$ws.Range("M46").Value = -33158
# Row 132: Lar Ingot
$ws.Range("H132").Value = 3269.6206
$ws.Range("I132").Value = 3247.1428
$ws.Range("K132").Value = 9741.428400000001
$ws.Range("M132").Value = -7211.428400000001
# Row 136: Pink Beryl
$ws.Range("H136").Value = 27583.75
$ws.Range("J136").Value = 27583.75
$ws.Range("L136").Value = 82751.25
$ws.Range("N136").Value = -87851.25

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Toad Leather
$ws.Range("H40").Value = 5743.0586
$ws.Range("I40").Value = 3967.5454
$ws.Range("J40").Value = 8998.166999999999
$ws.Range("K40").Value = 3967.5454
$ws.Range("L40").Value = 8998.166999999999
$ws.Range("M40").Value = -3831.5454
$ws.Range("N40").Value = -9270.166999999999
# Row 132: Silver Lobo Leather
$ws.Range("H132").Value = 254001
$ws.Range("I132").Value = 336501.34
$ws.Range("K132").Value = 1009504.02
$ws.Range("M132").Value = -1006974.02

$ws = $wb.Worksheets.Item("WVR")
# Row 113: Pixie Floss
$ws.Range("H113").Value = 2100.2188
$ws.Range("I113").Value = 2149.923
$ws.Range("K113").Value = 6449.768999999999
$ws.Range("M113").Value = -4279.768999999999
